# Replace the 4 placeholder <w:tab/> runs at the end of the "Name: " lines
# (in the DATA USER / DATA AUTHOR signature table) with real MERGEFIELD
# field codes for the user/author full name and affiliations, matching
# the "add data user and data author data in the end of agreement" edit.
#
# Word doesn't expose a COM call that edits individual runs' field codes
# directly, so we rebuild the whole paragraph's OOXML and swap it in via
# Range.InsertXML (the supported way to inject raw OOXML through this
# object model) on a Range that spans the complete target paragraph.

$d = $word.ActiveDocument

$wOpenXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Shared run properties used by every run on these "Name: ..." lines.
$rPrPlain      = '<w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/></w:rPr>'
$rPrUnderline  = '<w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:u w:val="single"/></w:rPr>'

function New-Run([string]$rPr, [string]$inner) {
    return "<w:r>$rPr$inner</w:r>"
}

# Builds the run sequence for one "Name: <<fullName>> <<affiliations>>"
# field pair: a MERGEFIELD built from begin/instrText/separate/end parts
# for the full name, a literal space, then a second MERGEFIELD (begin/
# instrText/separate/end) for the affiliations.
function New-NameFieldRuns([string]$fullNameField, [string]$affiliationsField) {
    $runs = @()
    $runs += New-Run $rPrUnderline '<w:fldChar w:fldCharType="begin"/>'
    $runs += New-Run $rPrUnderline '<w:instrText xml:space="preserve"> </w:instrText>'
    $runs += New-Run $rPrUnderline '<w:instrText>MERGEFIELD</w:instrText>'
    $runs += New-Run $rPrUnderline "<w:instrText xml:space=`"preserve`"> $fullNameField </w:instrText>"
    $runs += New-Run $rPrUnderline '<w:fldChar w:fldCharType="separate"/>'
    $runs += New-Run $rPrUnderline '<w:fldChar w:fldCharType="end"/>'
    $runs += New-Run $rPrUnderline '<w:t xml:space="preserve"> </w:t>'
    $runs += New-Run $rPrUnderline '<w:fldChar w:fldCharType="begin"/>'
    $runs += New-Run $rPrUnderline "<w:instrText xml:space=`"preserve`"> MERGEFIELD $affiliationsField </w:instrText>"
    $runs += New-Run $rPrUnderline '<w:fldChar w:fldCharType="separate"/>'
    $runs += New-Run $rPrUnderline '<w:fldChar w:fldCharType="end"/>'
    return [string]::Join('', $runs)
}

# Rebuilds the full "Name: " paragraph (same pPr / leading "Name: " run as
# the original template) with the tab runs replaced by the field runs.
function New-NameParagraphXml([string]$fullNameField, [string]$affiliationsField) {
    $pPr = '<w:pPr><w:spacing w:before="120"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:u w:val="single"/></w:rPr></w:pPr>'
    $nameLabelRun = New-Run $rPrPlain '<w:t xml:space="preserve">Name: </w:t>'
    $fieldRuns = New-NameFieldRuns $fullNameField $affiliationsField
    $paragraph = "<w:p w:rsidR=`"00135642`" w:rsidRDefault=`"00000000`">$pPr$nameLabelRun$fieldRuns</w:p>"

    return "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>" +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           "<pkg:xmlData><w:document $wOpenXmlNs><w:body>$paragraph</w:body></w:document></pkg:xmlData>" +
           '</pkg:part></pkg:package>'
}

# Locate the two "Name: <tab><tab><tab><tab>" paragraphs in document order:
# the first is in the DATA USER cell, the second in the DATA AUTHOR cell.
$targets = @()
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Name: `t`t`t`t`r") {
        $targets += $p
    }
}

if ($targets.Count -ne 2) {
    throw "Expected exactly 2 'Name:' placeholder paragraphs, found $($targets.Count)"
}

$userParagraph = $targets[0]
$userRange = $d.Range($userParagraph.Range.Start, $userParagraph.Range.End)
$userRange.InsertXML((New-NameParagraphXml "dataUserFullName" "dataUserAffiliations"))

$authorParagraph = $targets[1]
$authorRange = $d.Range($authorParagraph.Range.Start, $authorParagraph.Range.End)
$authorRange.InsertXML((New-NameParagraphXml "dataAuthorFullName" "dataAuthorAffiliations"))
